$p = $ppt.ActivePresentation

# Slide 5 contains a table whose style is switched from the custom
# "Table_0" style to the built-in table style.
$s = $p.Slides.Item(5)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $sh.Table.ApplyStyle("{53E84F77-6D7F-409A-BF1A-604D7B1B1ECF}")
    }
}
